# Applies the diff to the logbook "Submissions" table (the 2nd table in
# the document): column widths change, a double space is collapsed to a
# single space in one existing cell, and a previously-empty row is filled
# in with a new date / filename / description.

$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

# ---------------------------------------------------------------
# 1. Column widths: 988 -> 1375 dxa, 3543 -> 3430 dxa, 4485 -> 4211 dxa
#    (dxa / 20 = points). Setting the width on any cell in a column
#    resizes that column for every row in the table.
# ---------------------------------------------------------------
$t.Cell(1, 1).Width = 1375 / 20
$t.Cell(1, 2).Width = 3430 / 20
$t.Cell(1, 3).Width = 4211 / 20

# ---------------------------------------------------------------
# 2. Row 3 (07/07/2024): collapse the double space between "t2" and
#    "series" to a single space.
# ---------------------------------------------------------------
$row3 = $t.Cell(3, 3).Range
$row3.Find.Execute("t2  series", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "t2 series", 2)

# ---------------------------------------------------------------
# 3. Row 9 (previously completely empty): fill in date, filename and
#    description, matching the Calibri font used elsewhere in the table.
# ---------------------------------------------------------------
$t.Cell(9, 1).Range.Text = "22/08/2024"
$t.Cell(9, 1).Range.Paragraphs.Item(1).Range.Font.Name = "Calibri"

$t.Cell(9, 2).Range.Text = "Prostatexproject.ipynb"
$t.Cell(9, 2).Range.Paragraphs.Item(1).Range.Font.Name = "Calibri"

$t.Cell(9, 3).Range.Text = "Model updated"
$t.Cell(9, 3).Range.Paragraphs.Item(1).Range.Font.Name = "Calibri"
